$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date), Volumen, Precio promedio ponderado, and Precio $/Kg
# values between row 2 and row 4.

$ws.Range("D2").Value2 = 44874
$ws.Range("M2").Value  = 200
$ws.Range("P2").Value  = 7750
$ws.Range("S2").Value  = 7750

$ws.Range("D4").Value2 = 44923
$ws.Range("M4").Value  = 80
$ws.Range("P4").Value  = 7625
$ws.Range("S4").Value  = 7625
